$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: mark G7, H7, I7 as "done" (F7 already "done")
$ws.Range("G7").Value = "done"
$ws.Range("H7").Value = "done"
$ws.Range("I7").Value = "done"

# Row 8: E8 becomes a "Bad" styled empty cell (no border), G8/H8/I8 become "done"
$ws.Range("E8").Style = "Bad"
$ws.Range("G8").Value = "done"
$ws.Range("H8").Value = "done"
$ws.Range("I8").Value = "done"

# Row 9: E9, F9 become "Bad" styled empty cells (no border)
$ws.Range("E9:F9").Style = "Bad"

# Row 10: E10, F10, G10 become "Bad" styled empty cells (no border); I10 becomes "done"
$ws.Range("E10:G10").Style = "Bad"
$ws.Range("I10").Value = "done"

# Row 11: E11:H11 become "Bad" styled empty cells, keeping their bottom border
$ws.Range("E11:H11").Style = "Bad"
$ws.Range("E11:H11").Borders.Item(9).LineStyle = 1
$ws.Range("E11:H11").Borders.Item(9).Weight = 2

# Update sheet view: scroll so column B is the top-left visible column, and move the selection
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("F17").Select()
